# Edit script: merges the title textbox into the instructions textbox on
# slide 1 (new size/position + expanded copy), and adds a new slide 2 with
# a "legend" textbox describing the 1/2/3/4 upgrade keys.

$p = $ppt.ActivePresentation
$s1 = $p.Slides.Item(1)

# --- Slide 1: drop the big "<GAME TITLE>" textbox -------------------------
$titleBox = $s1.Shapes.Item(1)
$titleBox.Delete()

# What used to be "TextBox 2" (id 3) is now the first shape on the slide.
$body = $s1.Shapes.Item(1)
$bodyTr = $body.TextFrame.TextRange

$apos = [char]0x2019
$bodyTr.Text = "Press 1 to upgrade size`r" + `
  "`r" + `
  "Press 2 to upgrade speed`r" + `
  "`r" + `
  "Press 3 to upgrade income`r" + `
  "`r" + `
  "Press 4 to buy gun / ammo`r" + `
  "`r" + `
  "Collect    for 10 coins`r" + `
  "`r" + `
  "Don" + $apos + "t touch the void`r" + `
  "`r" + `
  "Use arrow keys to shoot`r" + `
  "`r" + `
  "Use WASD to move`r" + `
  "`r" + `
  "Upgrades require coins"

# Base look (keeps the inherited Courier New typeface, just resizes it).
$bodyTr.Font.Size = 42

# Per-word colour accents.
$bodyTr.Characters(20, 4).Font.Color.RGB = 15773696   # "size"   -> 00B0F0
$bodyTr.Characters(45, 5).Font.Color.RGB = 5287936    # "speed"  -> 00B050
$bodyTr.Characters(71, 6).Font.Color.RGB = 65535      # "income" -> FFFF00
$bodyTr.Characters(94, 3).Font.Color.RGB = 192        # "gun"    -> C00000
$bodyTr.Characters(100, 4).Font.Color.RGB = 192       # "ammo"   -> C00000

# Reposition / resize the box to its new, larger frame.
$body.Left = 188375 / 12700
$body.Top = 158980 / 12700
$body.Width = 15326984 / 12700
$body.Height = 11172253 / 12700

# --- Slide 2: new "legend" slide -------------------------------------------
$s2 = $p.Slides.Add(2, 12)
$legend = $s2.Shapes.AddTextbox(1, 0, 12134850 / 12700, 24377650 / 12700, 646331 / 12700)
$legendTr = $legend.TextFrame.TextRange
$legendTr.Text = "1 for Size           2 for Speed          3 for Income          4 for Gun"
$legendTr.Characters(7, 4).Font.Color.RGB = 12611584  # "Size" -> 0070C0

$legend.Fill.Visible = $false
$legend.TextFrame.WordWrap = $true
$legend.TextFrame.AutoSize = 1
$legend.Left = 0 / 12700
$legend.Top = 12134850 / 12700
$legend.Width = 24377650 / 12700
$legend.Height = 646331 / 12700
